$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.822.33"
$ws.Range("E2").Value = "  -1.53%  "
$ws.Range("D3").Value = "1.550.69"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'205.36"
$ws.Range("E5").Value = "  -1.28%  "
$ws.Range("D6").Value = "'0.482"
$ws.Range("E6").Value = "  -1.73%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'21.47"
$ws.Range("E8").Value = "  -3.40%  "
$ws.Range("D9").Value = "'0.246"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("D10").Value = "'0.0582"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "1.548.07"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("E14").Value = "  -2.48%  "
$ws.Range("E15").Value = "  -1.43%  "
$ws.Range("D16").Value = "26.792.71"
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "'61.07"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "'214.01"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "0.0₃0688"
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").Value = "'7.25"
$ws.Range("E20").Value = "  -1.30%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "'9.04"
$ws.Range("E23").Value = "  -3.78%  "
$ws.Range("E24").Value = "  -0.96%  "
$ws.Range("D25").Value = "'152.97"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").Value = "'6.52"
$ws.Range("E26").Value = "  -2.39%  "
$ws.Range("D27").Value = "'14.90"
$ws.Range("E27").Value = "  -0.34%  "
$ws.Range("E28").Value = "  -0.16%  "
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "'0.0462"
$ws.Range("E30").Value = "  -0.63%  "
$ws.Range("E31").Value = "  -3.20%  "
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "1.355.31"
$ws.Range("E33").Value = "  -3.64%  "
$ws.Range("D34").Value = "'2.91"
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("D35").Value = "'1.51"
$ws.Range("E35").Value = "  -3.37%  "
$ws.Range("D36").Value = "'2.26"
$ws.Range("E36").Value = "  -0.98%  "
$ws.Range("D37").Value = "'0.919"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("E38").Value = "  -2.04%  "
$ws.Range("E39").Value = "  +1.02%  "
$ws.Range("D40").Value = "'0.803"
$ws.Range("E40").Value = "  -1.90%  "
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("E42").Value = "  +4.29%  "
$ws.Range("D43").Value = "'0.992"
$ws.Range("E43").Value = "  -0.83%  "
$ws.Range("E44").Value = "  +0.18%  "
$ws.Range("D45").Value = "'1.77"
$ws.Range("E45").Value = "  -2.41%  "
$ws.Range("D46").Value = "'62.98"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("D48").Value = "1.683.73"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "'85.92"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  +2.54%  "
$ws.Range("D51").Value = "0.0₇0972"
$ws.Range("E51").Value = "  -1.85%  "
